# commit 2: new geniric classes and scripts
# Adds a "homepagelink" row to the "login" sheet (TestData.xlsx / sheet2 = login),
# widens the two data columns to fit the new (longer) content, and leaves the
# selection where Excel would land after typing the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # "login" sheet is the active tab in this workbook

# New row 4: key/value pair appended below the existing appurl/username/password rows
$ws.Range("A4").Value = "homepagelink"
$ws.Range("B4").Value = "http://localhost:8084/dashboard/welcome"

# Columns were re-fitted (bestFit) to the new, wider content
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 35.833333333333336

# Final selection left on C11 (as captured in the saved view state)
$ws.Range("C11").Select()
